# Updated cryptos list on Sun May 21 07:40:59 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row with
# the latest scraped figures, and reorders the ImmutableX / ARBITRUM rows
# (32 and 33 swap places, including their Coin name, Link and Price/Volume
# values).
#
# All cells on this sheet are plain text (coinranking.com renders prices
# like "27.324.64" / "1.834.40", which aren't valid numbers, and percentages
# like "  +1.59%  "). Whenever the new price text WOULD be auto-recognised
# by Excel as a genuine number (e.g. "1.011", "0.3690") we prefix it with a
# leading apostrophe, same as typing '1.011 into a cell by hand, so Excel
# keeps storing it as text instead of silently converting it to a number
# (and losing the trailing zero / dot-as-thousands-separator formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.324.64'
$ws.Range("E2").Value = '  +1.59%  '
$ws.Range("D3").Value = '1.834.40'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("D4").Value = '''1.011'
$ws.Range("E4").Value = '  +0.84%  '
$ws.Range("D5").Value = '''314.86'
$ws.Range("E5").Value = '  +1.96%  '
$ws.Range("D7").Value = '''0.4743'
$ws.Range("E7").Value = '  +1.91%  '
$ws.Range("D8").Value = '''0.3690'
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("D9").Value = '''0.07449'
$ws.Range("E9").Value = '  +1.23%  '
$ws.Range("D10").Value = '''0.8858'
$ws.Range("E10").Value = '  +2.05%  '
$ws.Range("D11").Value = '''20.47'
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("D12").Value = '1.870.74'
$ws.Range("E12").Value = '  +3.53%  '
$ws.Range("D13").Value = '''0.07347'
$ws.Range("E13").Value = '  +3.45%  '
$ws.Range("D14").Value = '''5.442'
$ws.Range("D15").Value = '''93.83'
$ws.Range("E15").Value = '  +3.03%  '
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("E17").Value = '  +0.69%  '
$ws.Range("D18").Value = '''0.000008795'
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("D19").Value = '''1.010'
$ws.Range("E19").Value = '  +0.89%  '
$ws.Range("E20").Value = '  +1.44%  '
$ws.Range("D21").Value = '27.512.11'
$ws.Range("E21").Value = '  +2.20%  '
$ws.Range("D22").Value = '''5.290'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("D24").Value = '2.088.20'
$ws.Range("E24").Value = '  +2.11%  '
$ws.Range("D25").Value = '''1.895'
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").Value = '''152.05'
$ws.Range("E26").Value = '  +0.76%  '
$ws.Range("D27").Value = '''18.66'
$ws.Range("E27").Value = '  +1.93%  '
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("D29").Value = '''5.232'
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("D30").Value = '''117.24'
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("E31").Value = '  +1.26%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '''0.7519'
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = '''1.177'
$ws.Range("E33").Value = '  +1.31%  '
$ws.Range("D34").Value = '''4.548'
$ws.Range("E34").Value = '  +1.65%  '
$ws.Range("D35").Value = '''2.949'
$ws.Range("E35").Value = '  +1.36%  '
$ws.Range("D36").Value = '''1.010'
$ws.Range("E36").Value = '  +0.82%  '
$ws.Range("E37").Value = '  +1.33%  '
$ws.Range("D38").Value = '''0.05348'
$ws.Range("E38").Value = '  +1.38%  '
$ws.Range("E39").Value = '  +0.71%  '
$ws.Range("D40").Value = '''2.956'
$ws.Range("E40").Value = '  -0.44%  '
$ws.Range("D41").Value = '''7.241'
$ws.Range("E41").Value = '  +0.97%  '
$ws.Range("D42").Value = '''2.385'
$ws.Range("E42").Value = '  +2.48%  '
$ws.Range("D43").Value = '''0.5312'
$ws.Range("E43").Value = '  +0.84%  '
$ws.Range("D44").Value = '''0.1660'
$ws.Range("E44").Value = '  +0.50%  '
$ws.Range("D45").Value = '''8.484'
$ws.Range("D46").Value = '''0.4926'
$ws.Range("E46").Value = '  +1.72%  '
$ws.Range("D47").Value = '''10.52'
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("D48").Value = '''105.12'
$ws.Range("E48").Value = '  +2.04%  '
$ws.Range("E49").Value = '  +0.96%  '
$ws.Range("D50").Value = '''1.674'
$ws.Range("E50").Value = '  +0.98%  '
$ws.Range("E51").Value = '  +0.09%  '
